$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarioA")
$ws.Name = "ScenarioA"
